# The commit swaps the DB connection source (the workbook used to read/write
# client CSV paths off a "G:\Database_automation\..." share; it now points at
# a local "D:\Database_automation\Database_automation\..." folder instead),
# and the previously-active "database" tab was left on "files" afterwards.

$wb = $excel.ActiveWorkbook

$dbSheet = $wb.Worksheets.Item("database")
$filesSheet = $wb.Worksheets.Item("files")

# Source_file_name / Target_file_name on the "files" sheet: same logical
# files (src stays src, target stays target), new connection path.
$filesSheet.Range("C2").Value = "D://Database_automation//Database_automation//files//client_data_src.csv"
$filesSheet.Range("F2").Value = "D://Database_automation//Database_automation//files//client_data_target.csv"

# "database" keeps its old selection (D12) but is no longer the active tab.
$dbSheet.Range("D12").Select() | Out-Null

# "files" becomes the active tab, selection moves to E12.
$filesSheet.Activate()
$filesSheet.Range("E12").Select() | Out-Null
